# Update countries & provincias Spain
#
# Refreshes the COVID-19 "Pais" sheet with a newer data pull:
#  - bumps the "last updated" banner in A1 to the 17:16 refresh
#  - updates case/recovery/death figures for several countries
#  - a handful of countries swapped ranking order (their whole data rows,
#    including the country name in column A, trade places) because the
#    newer pull re-sorted by total cases

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 17:16"

# Row 4 (country unchanged)
$ws.Range("B4").Value = 7104786
$ws.Range("C4").Value = 6849
$ws.Range("D4").Value = 4360074
$ws.Range("E4").Value = 2538913
$ws.Range("G4").Value = 328
$ws.Range("H4").Value = 205799

# Row 5 (country unchanged)
$ws.Range("B5").Value = 5669610
$ws.Range("C5").Value = 29114
$ws.Range("D5").Value = 4609704
$ws.Range("E5").Value = 969624
$ws.Range("G5").Value = 261
$ws.Range("H5").Value = 90282

# Row 19 -> Irak
$ws.Range("A19").Value = "Irak"
$ws.Range("B19").Value = 332635
$ws.Range("C19").Value = 5055
$ws.Range("D19").Value = 264988
$ws.Range("E19").Value = 58893
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = 8754

# Row 20 -> Arabia Saudita
$ws.Range("A20").Value = "Arabia Saudita"
$ws.Range("B20").Value = 331359
$ws.Range("C20").Value = 561
$ws.Range("D20").Value = 313786
$ws.Range("E20").Value = 13004
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = 4569

# Row 23 (country unchanged)
$ws.Range("B23").Value = 302537
$ws.Range("C23").Value = 1640
$ws.Range("D23").Value = 220665
$ws.Range("E23").Value = 46114
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = 35758

# Row 25 (country unchanged)
$ws.Range("B25").Value = 277877
$ws.Range("C25").Value = 701
$ws.Range("E25").Value = 20479

# Row 57 (country unchanged)
$ws.Range("D57").Value = 57291
$ws.Range("E57").Value = 321

# Row 63 -> Moldavia
$ws.Range("A63").Value = "Moldavia"
$ws.Range("B63").Value = 48232
$ws.Range("C63").Value = 786
$ws.Range("D63").Value = 36071
$ws.Range("E63").Value = 10917
$ws.Range("G63").Value = 14
$ws.Range("H63").Value = 1244

# Row 64 -> Armenia
$ws.Range("A64").Value = "Armenia"
$ws.Range("B64").Value = 47877
$ws.Range("C64").Value = 210
$ws.Range("D64").Value = 43026
$ws.Range("E64").Value = 3909
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 942

# Row 86 (country unchanged)
$ws.Range("B86").Value = 17049
$ws.Range("C86").Value = 182
$ws.Range("D86").Value = 14186
$ws.Range("E86").Value = 2153
$ws.Range("G86").Value = 5
$ws.Range("H86").Value = 710

# Row 94 (country unchanged)
$ws.Range("B94").Value = 12787
$ws.Range("C94").Value = 121
$ws.Range("D94").Value = 7139
$ws.Range("E94").Value = 5278
$ws.Range("G94").Value = 3
$ws.Range("H94").Value = 370

# Row 114 -> Jordania
$ws.Range("A114").Value = "Jordania"
$ws.Range("B114").Value = 6042
$ws.Range("C114").Value = 363
$ws.Range("D114").Value = 3812
$ws.Range("E114").Value = 2195
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 35

# Row 115 -> Malaui
$ws.Range("A115").Value = "Malaui"
$ws.Range("B115").Value = 5739
$ws.Range("D115").Value = 4065
$ws.Range("E115").Value = 1495
$ws.Range("H115").Value = 179

# Row 117 -> Jamaica
$ws.Range("A117").Value = "Jamaica"
$ws.Range("B117").Value = 5395
$ws.Range("C117").Value = 125
$ws.Range("D117").Value = 1444
$ws.Range("E117").Value = 3875
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 76

# Row 118 -> Cabo Verde
$ws.Range("A118").Value = "Cabo Verde"
$ws.Range("B118").Value = 5337
$ws.Range("D118").Value = 4742
$ws.Range("E118").Value = 543
$ws.Range("H118").Value = 52

# Row 119 -> Suazilandia
$ws.Range("A119").Value = "Suazilandia"
$ws.Range("B119").Value = 5307
$ws.Range("D119").Value = 4672
$ws.Range("E119").Value = 529
$ws.Range("H119").Value = 106

# Row 153 (country unchanged)
$ws.Range("B153").Value = 2029
$ws.Range("C153").Value = 1
$ws.Range("D153").Value = 1245
$ws.Range("E153").Value = 198

# Row 214 -> Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Row 215 -> Islas Malvinas
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
